$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 5 de Abril de 2020 a las 19:22'   # A1: 'Datos actualizados a 5 de Abril de 2020 a las 18:52' -> 'Datos actualizados a 5 de Abril de 2020 a las 19:22'
$ws.Cells.Item(4, 2).Value = 327848   # B4: 323953 -> 327848
$ws.Cells.Item(4, 3).Value = 16491   # C4: 12596 -> 16491
$ws.Cells.Item(4, 4).Value = 16700   # D4: 16598 -> 16700
$ws.Cells.Item(4, 5).Value = 301823   # E4: 298170 -> 301823
$ws.Cells.Item(4, 6).Value = 8519   # F4: 8474 -> 8519
$ws.Cells.Item(4, 7).Value = 873   # G4: 733 -> 873
$ws.Cells.Item(4, 8).Value = 9325   # H4: 9185 -> 9325
$ws.Cells.Item(12, 2).Value = 27069   # B12: 23934 -> 27069
$ws.Cells.Item(12, 3).Value = 3135   # C12: 0 -> 3135
$ws.Cells.Item(12, 4).Value = 1042   # D12: 786 -> 1042
$ws.Cells.Item(12, 5).Value = 25453   # E12: 22647 -> 25453
$ws.Cells.Item(12, 6).Value = 1381   # F12: 1311 -> 1381
$ws.Cells.Item(12, 7).Value = 73   # G12: 0 -> 73
$ws.Cells.Item(12, 8).Value = 574   # H12: 501 -> 574
$ws.Cells.Item(17, 2).Value = 12051   # B17: 11934 -> 12051
$ws.Cells.Item(17, 3).Value = 270   # C17: 153 -> 270
$ws.Cells.Item(17, 5).Value = 8849   # E17: 8732 -> 8849
$ws.Cells.Item(19, 2).Value = 10568   # B19: 10475 -> 10568
$ws.Cells.Item(19, 3).Value = 208   # C19: 115 -> 208
$ws.Cells.Item(19, 5).Value = 9986   # E19: 9900 -> 9986
$ws.Cells.Item(19, 7).Value = 10   # G19: 3 -> 10
$ws.Cells.Item(19, 8).Value = 455   # H19: 448 -> 455
$ws.Cells.Item(23, 1).Value = 'Noruega'   # A23: 'Australia' -> 'Noruega'
$ws.Cells.Item(23, 4).Value = 32   # D23: 2315 -> 32
$ws.Cells.Item(23, 5).Value = 5585   # E23: 3337 -> 5585
$ws.Cells.Item(23, 6).Value = 89   # F23: 91 -> 89
$ws.Cells.Item(23, 7).Value = 8   # G23: 5 -> 8
$ws.Cells.Item(23, 8).Value = 70   # H23: 35 -> 70
$ws.Cells.Item(24, 1).Value = 'Australia'   # A24: 'Noruega' -> 'Australia'
$ws.Cells.Item(24, 2).Value = 5687   # B24: 5686 -> 5687
$ws.Cells.Item(24, 3).Value = 137   # C24: 136 -> 137
$ws.Cells.Item(24, 4).Value = 2315   # D24: 32 -> 2315
$ws.Cells.Item(24, 5).Value = 3337   # E24: 5584 -> 3337
$ws.Cells.Item(24, 6).Value = 91   # F24: 89 -> 91
$ws.Cells.Item(24, 7).Value = 5   # G24: 8 -> 5
$ws.Cells.Item(24, 8).Value = 35   # H24: 70 -> 35
$ws.Cells.Item(28, 4).Value = 618   # D28: 528 -> 618
$ws.Cells.Item(28, 5).Value = 3819   # E28: 3909 -> 3819
$ws.Cells.Item(28, 6).Value = 307   # F28: 38 -> 307
$ws.Cells.Item(55, 6).Value = 46   # F55: 0 -> 46
$ws.Cells.Item(56, 4).Value = 320   # D56: 297 -> 320
$ws.Cells.Item(56, 5).Value = 983   # E56: 1006 -> 983
$ws.Cells.Item(68, 5).Value = 791   # E68: 792 -> 791
$ws.Cells.Item(68, 7).Value = 2   # G68: 1 -> 2
$ws.Cells.Item(68, 8).Value = 13   # H68: 12 -> 13
$ws.Cells.Item(73, 5).Value = 545   # E73: 547 -> 545
$ws.Cells.Item(73, 7).Value = 2   # G73: 0 -> 2
$ws.Cells.Item(73, 8).Value = 7   # H73: 5 -> 7
$ws.Cells.Item(91, 1).Value = 'Jordania'   # A91: 'Reunion' -> 'Jordania'
$ws.Cells.Item(91, 2).Value = 345   # B91: 344 -> 345
$ws.Cells.Item(91, 3).Value = 22   # C91: 10 -> 22
$ws.Cells.Item(91, 4).Value = 110   # D91: 40 -> 110
$ws.Cells.Item(91, 5).Value = 230   # E91: 304 -> 230
$ws.Cells.Item(91, 6).Value = 5   # F91: 4 -> 5
$ws.Cells.Item(91, 8).Value = 5   # H91: 0 -> 5
$ws.Cells.Item(92, 1).Value = 'Reunion'   # A92: 'Afganistan' -> 'Reunion'
$ws.Cells.Item(92, 2).Value = 344   # B92: 337 -> 344
$ws.Cells.Item(92, 3).Value = 10   # C92: 28 -> 10
$ws.Cells.Item(92, 4).Value = 40   # D92: 15 -> 40
$ws.Cells.Item(92, 5).Value = 304   # E92: 315 -> 304
$ws.Cells.Item(92, 6).Value = 4   # F92: 0 -> 4
$ws.Cells.Item(92, 8).Value = 0   # H92: 7 -> 0
$ws.Cells.Item(93, 1).Value = 'Afganistan'   # A93: 'Jordania' -> 'Afganistan'
$ws.Cells.Item(93, 2).Value = 337   # B93: 323 -> 337
$ws.Cells.Item(93, 3).Value = 28   # C93: 0 -> 28
$ws.Cells.Item(93, 4).Value = 15   # D93: 74 -> 15
$ws.Cells.Item(93, 5).Value = 315   # E93: 244 -> 315
$ws.Cells.Item(93, 6).Value = 0   # F93: 5 -> 0
$ws.Cells.Item(93, 8).Value = 7   # H93: 5 -> 7
$ws.Cells.Item(101, 2).Value = 234   # B101: 228 -> 234
$ws.Cells.Item(101, 3).Value = 17   # C101: 11 -> 17
$ws.Cells.Item(101, 5).Value = 210   # E101: 204 -> 210
$ws.Cells.Item(130, 1).Value = 'Monaco'   # A130: 'Madagascar' -> 'Monaco'
$ws.Cells.Item(130, 2).Value = 73   # B130: 72 -> 73
$ws.Cells.Item(130, 3).Value = 7   # C130: 2 -> 7
$ws.Cells.Item(130, 4).Value = 3   # D130: 2 -> 3
$ws.Cells.Item(130, 5).Value = 69   # E130: 70 -> 69
$ws.Cells.Item(130, 6).Value = 2   # F130: 6 -> 2
$ws.Cells.Item(130, 8).Value = 1   # H130: 0 -> 1
$ws.Cells.Item(131, 1).Value = 'Madagascar'   # A131: 'Monaco' -> 'Madagascar'
$ws.Cells.Item(131, 2).Value = 72   # B131: 66 -> 72
$ws.Cells.Item(131, 3).Value = 2   # C131: 0 -> 2
$ws.Cells.Item(131, 4).Value = 2   # D131: 3 -> 2
$ws.Cells.Item(131, 5).Value = 70   # E131: 62 -> 70
$ws.Cells.Item(131, 6).Value = 6   # F131: 2 -> 6
$ws.Cells.Item(131, 8).Value = 0   # H131: 1 -> 0
$ws.Cells.Item(156, 1).Value = 'Haiti'   # A156: 'Birmania' -> 'Haiti'
$ws.Cells.Item(156, 3).Value = 1   # C156: 0 -> 1
$ws.Cells.Item(156, 4).Value = 1   # D156: 0 -> 1
$ws.Cells.Item(156, 8).Value = 0   # H156: 1 -> 0
$ws.Cells.Item(157, 1).Value = 'Birmania'   # A157: 'Haiti' -> 'Birmania'
$ws.Cells.Item(157, 3).Value = 0   # C157: 1 -> 0
$ws.Cells.Item(157, 4).Value = 0   # D157: 1 -> 0
$ws.Cells.Item(157, 8).Value = 1   # H157: 0 -> 1
$ws.Cells.Item(182, 1).Value = 'Mozambique'   # A182: 'Surinam' -> 'Mozambique'
$ws.Cells.Item(182, 4).Value = 1   # D182: 0 -> 1
$ws.Cells.Item(182, 8).Value = 0   # H182: 1 -> 0
$ws.Cells.Item(183, 1).Value = 'Surinam'   # A183: 'Mozambique' -> 'Surinam'
$ws.Cells.Item(183, 4).Value = 0   # D183: 1 -> 0
$ws.Cells.Item(183, 8).Value = 1   # H183: 0 -> 1
$ws.Cells.Item(185, 1).Value = 'Republica del Chad'   # A185: 'Suazilandia' -> 'Republica del Chad'
$ws.Cells.Item(186, 1).Value = 'Suazilandia'   # A186: 'Republica del Chad' -> 'Suazilandia'
$ws.Cells.Item(196, 1).Value = 'Botsuana'   # A196: 'San Bartolome' -> 'Botsuana'
$ws.Cells.Item(196, 3).Value = 2   # C196: 0 -> 2
$ws.Cells.Item(196, 4).Value = 0   # D196: 1 -> 0
$ws.Cells.Item(196, 8).Value = 1   # H196: 0 -> 1
$ws.Cells.Item(197, 1).Value = 'San Bartolome'   # A197: 'Mauritania' -> 'San Bartolome'
$ws.Cells.Item(197, 4).Value = 1   # D197: 2 -> 1
$ws.Cells.Item(197, 5).Value = 5   # E197: 3 -> 5
$ws.Cells.Item(197, 8).Value = 0   # H197: 1 -> 0
$ws.Cells.Item(198, 1).Value = 'Mauritania'   # A198: 'Belice' -> 'Mauritania'
$ws.Cells.Item(198, 2).Value = 6   # B198: 5 -> 6
$ws.Cells.Item(198, 3).Value = 0   # C198: 1 -> 0
$ws.Cells.Item(198, 4).Value = 2   # D198: 0 -> 2
$ws.Cells.Item(198, 5).Value = 3   # E198: 5 -> 3
$ws.Cells.Item(198, 6).Value = 0   # F198: 1 -> 0
$ws.Cells.Item(198, 8).Value = 1   # H198: 0 -> 1
$ws.Cells.Item(199, 1).Value = 'Belice'   # A199: 'Nicaragua' -> 'Belice'
$ws.Cells.Item(199, 3).Value = 1   # C199: 0 -> 1
$ws.Cells.Item(199, 5).Value = 5   # E199: 4 -> 5
$ws.Cells.Item(199, 6).Value = 1   # F199: 0 -> 1
$ws.Cells.Item(199, 8).Value = 0   # H199: 1 -> 0
$ws.Cells.Item(200, 1).Value = 'Nicaragua'   # A200: 'Islas Turcas y Caicos' -> 'Nicaragua'
$ws.Cells.Item(200, 7).Value = 0   # G200: 1 -> 0
$ws.Cells.Item(201, 1).Value = 'Islas Turcas y Caicos'   # A201: 'Butan' -> 'Islas Turcas y Caicos'
$ws.Cells.Item(201, 4).Value = 0   # D201: 2 -> 0
$ws.Cells.Item(201, 5).Value = 4   # E201: 3 -> 4
$ws.Cells.Item(201, 7).Value = 1   # G201: 0 -> 1
$ws.Cells.Item(201, 8).Value = 1   # H201: 0 -> 1
$ws.Cells.Item(202, 1).Value = 'Butan'   # A202: 'Malaui' -> 'Butan'
$ws.Cells.Item(202, 2).Value = 5   # B202: 4 -> 5
$ws.Cells.Item(202, 4).Value = 2   # D202: 0 -> 2
$ws.Cells.Item(202, 5).Value = 3   # E202: 4 -> 3
$ws.Cells.Item(204, 1).Value = 'Malaui'   # A204: 'Botsuana' -> 'Malaui'
$ws.Cells.Item(204, 5).Value = 4   # E204: 3 -> 4
$ws.Cells.Item(204, 8).Value = 0   # H204: 1 -> 0
$ws.Cells.Item(206, 1).Value = 'Islas Virgenes Britanicas'   # A206: 'Anguila' -> 'Islas Virgenes Britanicas'
$ws.Cells.Item(207, 1).Value = 'Anguila'   # A207: 'Burundi' -> 'Anguila'
$ws.Cells.Item(208, 1).Value = 'Burundi'   # A208: 'Islas Virgenes Britanicas' -> 'Burundi'
